# edit.ps1 - Applies the DataPredict Libraries Licensing Agreement changes:
#  1. Remove proofErr around "TensorL" heading
#  2. Remove proofErr around "Robux" in Gross Revenue definition (merge runs)
#  3. Remove proofErr around "Robux" in "Revenue expressed in Robux must..." (merge runs)
#  4. Remove proofErr around "Robux" in "Any applicable exchange rates..." (merge runs)
#  5. Insert new 10.3 paragraph about disclosure of Software Libraries use
#  6. Add <w:lastRenderedPageBreak/> to the first run of the 11.1 paragraph

$d = $word.ActiveDocument

function Replace-ParagraphXml($paraIndex, $innerRunsXml) {
    $p = $d.Paragraphs.Item($paraIndex)
    $full = $p.Range
    $target = $d.Range($full.Start, $full.End - 1)
    $frag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + "`n" `
        + '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' `
        + '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' `
        + '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' `
        + $innerRunsXml `
        + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $target.InsertXML($frag)
}

function Find-ParagraphIndexStartingWith($prefix) {
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $t = $d.Paragraphs.Item($i).Range.Text
        if ($t.StartsWith($prefix)) {
            return $i
        }
    }
    return -1
}

$apos = [char]0x2019
$ldq = [char]0x201C
$rdq = [char]0x201D

# ---------------------------------------------------------------------------
# 1. "TensorL:" heading - drop proofErr spellStart/spellEnd around "TensorL"
# ---------------------------------------------------------------------------
$idx1 = Find-ParagraphIndexStartingWith("TensorL:")
if ($idx1 -eq -1) { throw "Could not find TensorL: paragraph" }
$xml1 = '<w:r><w:t>TensorL</w:t></w:r><w:r><w:t>:</w:t></w:r>'
Replace-ParagraphXml $idx1 $xml1

# ---------------------------------------------------------------------------
# 2. Gross Revenue definition paragraph - merge the 3 runs around "Robux"
# ---------------------------------------------------------------------------
$idx2 = Find-ParagraphIndexStartingWith($ldq + "Gross Revenue" + $rdq + " refers")
if ($idx2 -eq -1) { throw "Could not find Gross Revenue definition paragraph" }
$text2 = $ldq + "Gross Revenue" + $rdq + " refers to the total monthly revenue earned from the Licensed Experience, including revenue derived from Robux transactions, Developer Products, Game Passes, in-game purchases, or any monetization mechanisms tied to the use of the Software Libraries."
$xml2 = '<w:r><w:t>' + $text2 + '</w:t></w:r>'
Replace-ParagraphXml $idx2 $xml2

# ---------------------------------------------------------------------------
# 3. "Gross Revenue expressed in Robux must..." paragraph - merge 3 runs
# ---------------------------------------------------------------------------
$idx3 = Find-ParagraphIndexStartingWith("Gross Revenue expressed in Robux must")
if ($idx3 -eq -1) { throw "Could not find Revenue expressed in Robux paragraph" }
$xml3 = '<w:r><w:t xml:space="preserve">Gross </w:t></w:r>' `
    + '<w:r><w:t xml:space="preserve">Revenue expressed in Robux </w:t></w:r>' `
    + '<w:r><w:t>must</w:t></w:r>' `
    + '<w:r><w:t xml:space="preserve"> be converted to its USD equivalent using the current or most recent Roblox Developer Exchange Rate as published by Roblox Corporation</w:t></w:r>' `
    + '<w:r><w:t xml:space="preserve"> </w:t></w:r>' `
    + '<w:r><w:t>before</w:t></w:r>' `
    + '<w:r><w:t xml:space="preserve"> </w:t></w:r>' `
    + '<w:r><w:t>calculating the license fee</w:t></w:r>' `
    + '<w:r><w:t>.</w:t></w:r>'
Replace-ParagraphXml $idx3 $xml3

# ---------------------------------------------------------------------------
# 4. "Any applicable exchange rates (e.g., Robux-to-USD conversions)," - merge
# ---------------------------------------------------------------------------
$idx4 = Find-ParagraphIndexStartingWith("Any applicable exchange rates")
if ($idx4 -eq -1) { throw "Could not find Any applicable exchange rates paragraph" }
$text4 = "Any applicable exchange rates (e.g., Robux-to-USD conversions),"
$xml4 = '<w:r><w:t>' + $text4 + '</w:t></w:r>'
Replace-ParagraphXml $idx4 $xml4

# ---------------------------------------------------------------------------
# 5. Insert new 10.3 paragraph after the 10.2 paragraph
# ---------------------------------------------------------------------------
$idx5 = Find-ParagraphIndexStartingWith("10.2. Licensee agrees that any gross revenue")
if ($idx5 -eq -1) { throw "Could not find 10.2 paragraph" }
$p5 = $d.Paragraphs.Item($idx5)
$full5 = $p5.Range
$insertPoint5 = $d.Range($full5.End - 1, $full5.End - 1)
$newParaXml = '<w:r><w:t xml:space="preserve">10.3. Licensee agrees to disclose the use of Software Libraries </w:t></w:r>' `
    + '<w:r><w:t xml:space="preserve">and the Licensor' + $apos + 's identity </w:t></w:r>' `
    + '<w:r><w:t>in all documentation</w:t></w:r>' `
    + '<w:r><w:t xml:space="preserve"> when applicable</w:t></w:r>' `
    + '<w:r><w:t>.</w:t></w:r>'
$frag5 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + "`n" `
    + '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' `
    + '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' `
    + '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' `
    + $newParaXml `
    + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertPoint5.InsertXML($frag5)

# ---------------------------------------------------------------------------
# 6. Add <w:lastRenderedPageBreak/> before "1" in the 11.1 paragraph
# ---------------------------------------------------------------------------
$idx6 = Find-ParagraphIndexStartingWith("11.1. This Agreement constitutes the entire agreement")
if ($idx6 -eq -1) { throw "Could not find 11.1 paragraph" }
$bodyText6 = "This Agreement constitutes the entire agreement between the parties with respect to the subject matter hereof and supersedes all prior and contemporaneous agreements, understandings, and representations, whether oral or written."
$xml6 = '<w:r><w:lastRenderedPageBreak/><w:t>1</w:t></w:r>' `
    + '<w:r><w:t>1</w:t></w:r>' `
    + '<w:r><w:t xml:space="preserve">.1. ' + $bodyText6 + '</w:t></w:r>'
Replace-ParagraphXml $idx6 $xml6

Write-Output "All edits applied."
